$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")
$summary = $wb.Worksheets.Item("Summary")

# New order-item rows appended after the existing last row (61).
# Numeric-looking values are written with a leading apostrophe so they are
# stored as text (matching the sheet's existing "number stored as text" data).

$ws.Range("C61").Value = "4_阳光粉洋桔梗_Sunshine Pink Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"
$ws.Range("F61").Value = "'20"

$ws.Range("C62").Value = "157_流沙_Quicksand_Rosa rugosa Thunb._20stems"
$ws.Range("F62").Value = "'7"

$ws.Range("C63").Value = "137_凯瑟琳_Catherine_Rosa rugosa Thunb._20stems"
$ws.Range("F63").Value = "'1"

$ws.Range("C64").Value = "277_草莓杏仁饼_undefined_Rosa rugosa Thunb._10stems"
$ws.Range("F64").Value = "'15"

$ws.Range("C65").Value = "238_苏菲宝贝_undefined_Rosa rugosa Thunb._10stems"
$ws.Range("F65").Value = "'5"

$ws.Range("C66").Value = "315_尤加利叶圆叶_Eucalyptus Populus_undefined_1bunch"
$ws.Range("F66").Value = "'5"

$ws.Range("C67").Value = "317_尤加利叶细叶_Eucalyptus Parvifolia_undefined_1bunch"
$ws.Range("F67").Value = "'5"

$ws.Range("C68").Value = "316_尤加利叶大叶_Eucalyptus Cinerea_undefined_1bunch"
$ws.Range("F68").Value = "'5"

$ws.Range("C69").Value = "578_腊梅粉_wax pink_undefined_1bunch"
$ws.Range("F69").Value = "'10"

$ws.Range("A70").Value = "'3"
$ws.Range("C70").Value = "327_文竹_asparagus fern_undefined_1bunch"
$ws.Range("F70").Value = "'5"

$ws.Range("C71").Value = "522_山归来绿_Smilax china_undefined_1bunch"
$ws.Range("F71").Value = "'3"

$ws.Range("C72").Value = "448_吊米 绿_hanging amaranthus`ngreen_undefined_1bunch"
$ws.Range("F72").Value = "'5"

$ws.Range("C73").Value = "557_山茶花_camellia_undefined_1bunch"
$ws.Range("F73").Value = "'3"

$ws.Range("C74").Value = "431_小米果_undefined_undefined_1bunch"
$ws.Range("F74").Value = "'5"

$ws.Range("C75").Value = "816_山里红_undefined_undefined_1bunch"
$ws.Range("F75").Value = "'10"

$ws.Range("C76").Value = "816_山里红_undefined_undefined_1bunch"
$ws.Range("F76").Value = "'8"

$ws.Range("C77").Value = "359_蓝梦叶_undefined_undefined_1bunch"
$ws.Range("F77").Value = "'5"

$ws.Range("C78").Value = "349_千层金绿_Melaleuca bracteata`n（dyed orange）_Melaleuca bracteata F.Muell._1bunch"
$ws.Range("F78").Value = "'5"

$ws.Range("C79").Value = "320_雪柳花_Spiraea flower white_undefined_1bunch"
$ws.Range("F79").Value = "'5"

$ws.Range("C80").Value = "100_绣球单瓣白_Hydrangea White S_Hydrangea L._1stem"
$ws.Range("F80").Value = "'15"

$ws.Range("C81").Value = "348_万年青_undefined_undefined_1bunch"

# Summary sheet: G2 is the concatenation of "0" + all Orders!F values (rows 2..81) + "0".
# Leading apostrophe forces text storage so this long digit string isn't
# collapsed into floating-point (scientific notation) and loses precision.
$summary.Range("G2").Value = "'0161112101016141010301010103050401052010101015151020103555151510201510101015151051010101540910102610101317810520711555551053535108555150"
